$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 5 (the "Reflection Pool v3 (No PadNFlute)" row no longer exists)
$ws.Rows.Item(5).Delete()

# Force Track Number column (L) to text so numeric-looking values keep their string type
$ws.Range("L2:L4").NumberFormat = "@"

# Update cell values for rows 2-4 to reflect corrected / new track metadata
$ws.Range("A2").Value = "DLM - Another Reminder v1 (Full).wav"
$ws.Range("AB2").Value = "US-RRD-20-00002"
$ws.Range("AK2").Value = "Hal"
$ws.Range("AL2").Value = "Stephens"
$ws.Range("AS2").Value = ""
$ws.Range("AZ2").Value = ""
$ws.Range("BH2").Value = ""
$ws.Range("BO2").Value = ""
$ws.Range("BW2").Value = ""
$ws.Range("CD2").Value = ""
$ws.Range("CL2").Value = ""
$ws.Range("CS2").Value = ""
$ws.Range("DC2").Value = "R42"
$ws.Range("E2").Value = "Another Reminder v4 (No Risers)"
$ws.Range("F2").Value = "DLM - Another Reminder v4 (No Risers).wav"
$ws.Range("H2").Value = "Jazz, Smooth Jazz Vol.42"
$ws.Range("I2").Value = "DLM-BI-176-R42"
$ws.Range("J2").Value = "Jazz (Smooth)"
$ws.Range("K2").Value = "Jazz (Smooth)"
$ws.Range("L2").Value = "2020422"
$ws.Range("M2").Value = "Hal  Stephens "
$ws.Range("N2").Value = "Hal  Stephens "
$ws.Range("P2").Value = "Jazz (Smooth)"
$ws.Range("T2").Value = "08-20-2020"
$ws.Range("W2").Value = "Smooth Jazz"
$ws.Range("AD3").Value = ""
$ws.Range("AK3").Value = ""
$ws.Range("AS3").Value = ""
$ws.Range("AZ3").Value = ""
$ws.Range("BH3").Value = ""
$ws.Range("BO3").Value = ""
$ws.Range("BW3").Value = ""
$ws.Range("CD3").Value = ""
$ws.Range("CL3").Value = ""
$ws.Range("CS3").Value = ""
$ws.Range("DC3").Value = "R42"
$ws.Range("E3").Value = "Nothing In This World v1 (Full)"
$ws.Range("F3").Value = "DLM - Nothing In This World v1 (Full).wav"
$ws.Range("H3").Value = "N/A, N/A Vol.42"
$ws.Range("I3").Value = "DLM-BI-NaN-R42"
$ws.Range("L3").Value = "2020423"
$ws.Range("T3").Value = "08-20-2020"
$ws.Range("AB4").Value = "US-RRD-20-00001"
$ws.Range("AD4").Value = ""
$ws.Range("AK4").Value = ""
$ws.Range("AS4").Value = ""
$ws.Range("AZ4").Value = ""
$ws.Range("BH4").Value = ""
$ws.Range("BO4").Value = ""
$ws.Range("BW4").Value = ""
$ws.Range("CD4").Value = ""
$ws.Range("CL4").Value = ""
$ws.Range("CS4").Value = ""
$ws.Range("DC4").Value = "R42"
$ws.Range("E4").Value = "Another Reminder v1 (Full)"
$ws.Range("F4").Value = "DLM - Another Reminder v1 (Full).wav"
$ws.Range("H4").Value = "Jazz, Smooth Jazz Vol.42"
$ws.Range("I4").Value = "DLM-BI-176-R42"
$ws.Range("L4").Value = "2020421"
$ws.Range("T4").Value = "08-20-2020"
